$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "/home/test/Desktop/Sleep/code-revised/data/SC4001E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4002E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4011E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4012E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4021E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4022E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4031E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4032E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4041E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4042E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4051E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4052E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4061E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4062E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4071E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4072E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4081E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4082E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4091E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4092E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4101E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4102E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4111E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4112E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4121E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4122E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4131E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4141E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4142E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4151E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4152E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4161E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4162E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4171E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4172E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4181E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4182E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4191E0.npz",
    "/home/test/Desktop/Sleep/code-revised/data/SC4192E0.npz"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws.Range("G11").Select() | Out-Null

